# Applies the RPA datasets push 2024-01-09 edit to the IPO underwriting table.
# The source row for NH / 캡스톤파트너스 (2023-11-06) was dropped from the feed,
# and the remaining rows were reordered within each underwriter block as part of
# the refreshed data pull. We trim the now-unused last row, then rewrite every
# data cell (rows 2-27) explicitly so the sheet matches the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One fewer data row in the refreshed export -- drop the now-unused last row.
$ws.Rows(28).Delete()

# Row 2: IBK / IBKS제23호스팩
$ws.Cells.Item(2, 1).Value = 'IBK'
$ws.Cells.Item(2, 2).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(2, 2).PasteSpecial(-4163)
$ws.Cells.Item(2, 3).Value = 'IBKS제23호스팩'
$ws.Cells.Item(2, 4).Value = 'IBK'
$ws.Cells.Item(2, 5).Value = 'IBK'
$ws.Cells.Item(2, 6).Formula = '=""&"2023-12-15"'
$ws.Cells.Item(2, 6).Copy()
$ws.Cells.Item(2, 6).PasteSpecial(-4163)
$ws.Cells.Item(2, 7).Formula = '=""&"2023-12-22"'
$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item(2, 7).PasteSpecial(-4163)
$ws.Cells.Item(2, 8).Value = 8000
$ws.Cells.Item(2, 9).Value = 4000000
$ws.Cells.Item(2, 10).Value = 2000
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 100

# Row 3: KB / DS단석
$ws.Cells.Item(3, 1).Value = 'KB'
$ws.Cells.Item(3, 2).Formula = '=""&"2023-12-14"'
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(3, 2).PasteSpecial(-4163)
$ws.Cells.Item(3, 3).Value = 'DS단석'
$ws.Cells.Item(3, 4).Value = 'KB, NH'
$ws.Cells.Item(3, 5).Value = 'KB, NH'
$ws.Cells.Item(3, 6).Formula = '=""&"2023-12-19"'
$ws.Cells.Item(3, 6).Copy()
$ws.Cells.Item(3, 6).PasteSpecial(-4163)
$ws.Cells.Item(3, 7).Formula = '=""&"2023-12-22"'
$ws.Cells.Item(3, 7).Copy()
$ws.Cells.Item(3, 7).PasteSpecial(-4163)
$ws.Cells.Item(3, 8).Value = 79300
$ws.Cells.Item(3, 9).Value = 1220000
$ws.Cells.Item(3, 10).Value = 100000
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 65

# Row 4: KB / LS머트리얼즈
$ws.Cells.Item(4, 1).Value = 'KB'
$ws.Cells.Item(4, 2).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(4, 2).Copy()
$ws.Cells.Item(4, 2).PasteSpecial(-4163)
$ws.Cells.Item(4, 3).Value = 'LS머트리얼즈'
$ws.Cells.Item(4, 4).Value = '키움, KB'
$ws.Cells.Item(4, 5).Value = '키움, KB, 이베스트, 하이, NH'
$ws.Cells.Item(4, 6).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(4, 6).Copy()
$ws.Cells.Item(4, 6).PasteSpecial(-4163)
$ws.Cells.Item(4, 7).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(4, 7).Copy()
$ws.Cells.Item(4, 7).PasteSpecial(-4163)
$ws.Cells.Item(4, 8).Value = 36196.872
$ws.Cells.Item(4, 9).Value = 14625000
$ws.Cells.Item(4, 10).Value = 6000
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 41.25

# Row 5: KB / 에코아이
$ws.Cells.Item(5, 1).Value = 'KB'
$ws.Cells.Item(5, 2).Formula = '=""&"2023-11-10"'
$ws.Cells.Item(5, 2).Copy()
$ws.Cells.Item(5, 2).PasteSpecial(-4163)
$ws.Cells.Item(5, 3).Value = '에코아이'
$ws.Cells.Item(5, 4).Value = 'KB'
$ws.Cells.Item(5, 5).Value = 'KB'
$ws.Cells.Item(5, 6).Formula = '=""&"2023-11-15"'
$ws.Cells.Item(5, 6).Copy()
$ws.Cells.Item(5, 6).PasteSpecial(-4163)
$ws.Cells.Item(5, 7).Formula = '=""&"2023-11-21"'
$ws.Cells.Item(5, 7).Copy()
$ws.Cells.Item(5, 7).PasteSpecial(-4163)
$ws.Cells.Item(5, 8).Value = 72141.3
$ws.Cells.Item(5, 9).Value = 2079000
$ws.Cells.Item(5, 10).Value = 34700
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 100

# Row 6: NH / DS단석
$ws.Cells.Item(6, 1).Value = 'NH'
$ws.Cells.Item(6, 2).Formula = '=""&"2023-12-14"'
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4163)
$ws.Cells.Item(6, 3).Value = 'DS단석'
$ws.Cells.Item(6, 4).Value = 'KB, NH'
$ws.Cells.Item(6, 5).Value = 'KB, NH'
$ws.Cells.Item(6, 6).Formula = '=""&"2023-12-19"'
$ws.Cells.Item(6, 6).Copy()
$ws.Cells.Item(6, 6).PasteSpecial(-4163)
$ws.Cells.Item(6, 7).Formula = '=""&"2023-12-22"'
$ws.Cells.Item(6, 7).Copy()
$ws.Cells.Item(6, 7).PasteSpecial(-4163)
$ws.Cells.Item(6, 8).Value = 42700
$ws.Cells.Item(6, 9).Value = 1220000
$ws.Cells.Item(6, 10).Value = 100000
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 35

# Row 7: NH / 동인기연
$ws.Cells.Item(7, 1).Value = 'NH'
$ws.Cells.Item(7, 2).Formula = '=""&"2023-11-09"'
$ws.Cells.Item(7, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)
$ws.Cells.Item(7, 3).Value = '동인기연'
$ws.Cells.Item(7, 4).Value = 'NH'
$ws.Cells.Item(7, 5).Value = 'NH'
$ws.Cells.Item(7, 6).Formula = '=""&"2023-11-14"'
$ws.Cells.Item(7, 6).Copy()
$ws.Cells.Item(7, 6).PasteSpecial(-4163)
$ws.Cells.Item(7, 7).Formula = '=""&"2023-11-22"'
$ws.Cells.Item(7, 7).Copy()
$ws.Cells.Item(7, 7).PasteSpecial(-4163)
$ws.Cells.Item(7, 8).Value = 44112
$ws.Cells.Item(7, 9).Value = 1470400
$ws.Cells.Item(7, 10).Value = 30000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 100

# Row 8: NH / LS머트리얼즈
$ws.Cells.Item(8, 1).Value = 'NH'
$ws.Cells.Item(8, 2).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(8, 2).PasteSpecial(-4163)
$ws.Cells.Item(8, 3).Value = 'LS머트리얼즈'
$ws.Cells.Item(8, 4).Value = '키움, KB'
$ws.Cells.Item(8, 5).Value = '키움, KB, 이베스트, 하이, NH'
$ws.Cells.Item(8, 6).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(8, 6).Copy()
$ws.Cells.Item(8, 6).PasteSpecial(-4163)
$ws.Cells.Item(8, 7).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(8, 7).Copy()
$ws.Cells.Item(8, 7).PasteSpecial(-4163)
$ws.Cells.Item(8, 8).Value = 4387.5
$ws.Cells.Item(8, 9).Value = 14625000
$ws.Cells.Item(8, 10).Value = 6000
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 5

# Row 9: NH / 그린리소스
$ws.Cells.Item(9, 1).Value = 'NH'
$ws.Cells.Item(9, 2).Formula = '=""&"2023-11-13"'
$ws.Cells.Item(9, 2).Copy()
$ws.Cells.Item(9, 2).PasteSpecial(-4163)
$ws.Cells.Item(9, 3).Value = '그린리소스'
$ws.Cells.Item(9, 4).Value = 'NH'
$ws.Cells.Item(9, 5).Value = 'NH'
$ws.Cells.Item(9, 6).Formula = '=""&"2023-11-16"'
$ws.Cells.Item(9, 6).Copy()
$ws.Cells.Item(9, 6).PasteSpecial(-4163)
$ws.Cells.Item(9, 7).Formula = '=""&"2023-11-24"'
$ws.Cells.Item(9, 7).Copy()
$ws.Cells.Item(9, 7).PasteSpecial(-4163)
$ws.Cells.Item(9, 8).Value = 27880
$ws.Cells.Item(9, 9).Value = 1640000
$ws.Cells.Item(9, 10).Value = 17000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 100

# Row 10: NH / 에코프로머티
$ws.Cells.Item(10, 1).Value = 'NH'
$ws.Cells.Item(10, 2).Formula = '=""&"2023-11-08"'
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(10, 2).PasteSpecial(-4163)
$ws.Cells.Item(10, 3).Value = '에코프로머티'
$ws.Cells.Item(10, 4).Value = '미래'
$ws.Cells.Item(10, 5).Value = '미래, NH, 하이'
$ws.Cells.Item(10, 6).Formula = '=""&"2023-11-13"'
$ws.Cells.Item(10, 6).Copy()
$ws.Cells.Item(10, 6).PasteSpecial(-4163)
$ws.Cells.Item(10, 7).Formula = '=""&"2023-11-17"'
$ws.Cells.Item(10, 7).Copy()
$ws.Cells.Item(10, 7).PasteSpecial(-4163)
$ws.Cells.Item(10, 8).Value = 121994.4706
$ws.Cells.Item(10, 9).Value = 11580800
$ws.Cells.Item(10, 10).Value = 36200
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 29.1

# Row 11: NH / 엔에이치스팩30호
$ws.Cells.Item(11, 1).Value = 'NH'
$ws.Cells.Item(11, 2).Formula = '=""&"2023-11-21"'
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4163)
$ws.Cells.Item(11, 3).Value = '엔에이치스팩30호'
$ws.Cells.Item(11, 4).Value = 'NH'
$ws.Cells.Item(11, 5).Value = 'NH'
$ws.Cells.Item(11, 6).Formula = '=""&"2023-11-24"'
$ws.Cells.Item(11, 6).Copy()
$ws.Cells.Item(11, 6).PasteSpecial(-4163)
$ws.Cells.Item(11, 7).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(11, 7).Copy()
$ws.Cells.Item(11, 7).PasteSpecial(-4163)
$ws.Cells.Item(11, 8).Value = 16000
$ws.Cells.Item(11, 9).Value = 8000000
$ws.Cells.Item(11, 10).Value = 2000
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 100

# Row 12: 교보 / 교보15호스팩
$ws.Cells.Item(12, 1).Value = '교보'
$ws.Cells.Item(12, 2).Formula = '=""&"2023-11-23"'
$ws.Cells.Item(12, 2).Copy()
$ws.Cells.Item(12, 2).PasteSpecial(-4163)
$ws.Cells.Item(12, 3).Value = '교보15호스팩'
$ws.Cells.Item(12, 4).Value = '교보'
$ws.Cells.Item(12, 5).Value = '교보'
$ws.Cells.Item(12, 6).Formula = '=""&"2023-11-28"'
$ws.Cells.Item(12, 6).Copy()
$ws.Cells.Item(12, 6).PasteSpecial(-4163)
$ws.Cells.Item(12, 7).Formula = '=""&"2023-12-05"'
$ws.Cells.Item(12, 7).Copy()
$ws.Cells.Item(12, 7).PasteSpecial(-4163)
$ws.Cells.Item(12, 8).Value = 7000
$ws.Cells.Item(12, 9).Value = 3500000
$ws.Cells.Item(12, 10).Value = 2000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 100

# Row 13: 대신 / 한선엔지니어링
$ws.Cells.Item(13, 1).Value = '대신'
$ws.Cells.Item(13, 2).Formula = '=""&"2023-11-13"'
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4163)
$ws.Cells.Item(13, 3).Value = '한선엔지니어링'
$ws.Cells.Item(13, 4).Value = '대신'
$ws.Cells.Item(13, 5).Value = '대신'
$ws.Cells.Item(13, 6).Formula = '=""&"2023-11-16"'
$ws.Cells.Item(13, 6).Copy()
$ws.Cells.Item(13, 6).PasteSpecial(-4163)
$ws.Cells.Item(13, 7).Formula = '=""&"2023-11-24"'
$ws.Cells.Item(13, 7).Copy()
$ws.Cells.Item(13, 7).PasteSpecial(-4163)
$ws.Cells.Item(13, 8).Value = 29750
$ws.Cells.Item(13, 9).Value = 4250000
$ws.Cells.Item(13, 10).Value = 7000
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 100

# Row 14: 미래 / 에코프로머티
$ws.Cells.Item(14, 1).Value = '미래'
$ws.Cells.Item(14, 2).Formula = '=""&"2023-11-08"'
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)
$ws.Cells.Item(14, 3).Value = '에코프로머티'
$ws.Cells.Item(14, 4).Value = '미래'
$ws.Cells.Item(14, 5).Value = '미래, NH, 하이'
$ws.Cells.Item(14, 6).Formula = '=""&"2023-11-13"'
$ws.Cells.Item(14, 6).Copy()
$ws.Cells.Item(14, 6).PasteSpecial(-4163)
$ws.Cells.Item(14, 7).Formula = '=""&"2023-11-17"'
$ws.Cells.Item(14, 7).Copy()
$ws.Cells.Item(14, 7).PasteSpecial(-4163)
$ws.Cells.Item(14, 8).Value = 284653.7406
$ws.Cells.Item(14, 9).Value = 11580800
$ws.Cells.Item(14, 10).Value = 36200
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 67.9

# Row 15: 미래 / 에이에스텍
$ws.Cells.Item(15, 1).Value = '미래'
$ws.Cells.Item(15, 2).Formula = '=""&"2023-11-16"'
$ws.Cells.Item(15, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4163)
$ws.Cells.Item(15, 3).Value = '에이에스텍'
$ws.Cells.Item(15, 4).Value = '미래'
$ws.Cells.Item(15, 5).Value = '미래'
$ws.Cells.Item(15, 6).Formula = '=""&"2023-11-21"'
$ws.Cells.Item(15, 6).Copy()
$ws.Cells.Item(15, 6).PasteSpecial(-4163)
$ws.Cells.Item(15, 7).Formula = '=""&"2023-11-28"'
$ws.Cells.Item(15, 7).Copy()
$ws.Cells.Item(15, 7).PasteSpecial(-4163)
$ws.Cells.Item(15, 8).Value = 39396
$ws.Cells.Item(15, 9).Value = 1407000
$ws.Cells.Item(15, 10).Value = 28000
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 100

# Row 16: 삼성 / 삼성스팩9호
$ws.Cells.Item(16, 1).Value = '삼성'
$ws.Cells.Item(16, 2).Formula = '=""&"2023-11-23"'
$ws.Cells.Item(16, 2).Copy()
$ws.Cells.Item(16, 2).PasteSpecial(-4163)
$ws.Cells.Item(16, 3).Value = '삼성스팩9호'
$ws.Cells.Item(16, 4).Value = '삼성'
$ws.Cells.Item(16, 5).Value = '삼성'
$ws.Cells.Item(16, 6).Formula = '=""&"2023-11-28"'
$ws.Cells.Item(16, 6).Copy()
$ws.Cells.Item(16, 6).PasteSpecial(-4163)
$ws.Cells.Item(16, 7).Formula = '=""&"2023-12-04"'
$ws.Cells.Item(16, 7).Copy()
$ws.Cells.Item(16, 7).PasteSpecial(-4163)
$ws.Cells.Item(16, 8).Value = 20000
$ws.Cells.Item(16, 9).Value = 10000000
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 100

# Row 17: 신영 / 케이엔에스
$ws.Cells.Item(17, 1).Value = '신영'
$ws.Cells.Item(17, 2).Formula = '=""&"2023-11-27"'
$ws.Cells.Item(17, 2).Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4163)
$ws.Cells.Item(17, 3).Value = '케이엔에스'
$ws.Cells.Item(17, 4).Value = '신영'
$ws.Cells.Item(17, 5).Value = '신영'
$ws.Cells.Item(17, 6).Formula = '=""&"2023-11-30"'
$ws.Cells.Item(17, 6).Copy()
$ws.Cells.Item(17, 6).PasteSpecial(-4163)
$ws.Cells.Item(17, 7).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(17, 7).Copy()
$ws.Cells.Item(17, 7).PasteSpecial(-4163)
$ws.Cells.Item(17, 8).Value = 17250
$ws.Cells.Item(17, 9).Value = 750000
$ws.Cells.Item(17, 10).Value = 23000
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 100

# Row 18: 유안타 / 와이바이오로직스
$ws.Cells.Item(18, 1).Value = '유안타'
$ws.Cells.Item(18, 2).Formula = '=""&"2023-11-23"'
$ws.Cells.Item(18, 2).Copy()
$ws.Cells.Item(18, 2).PasteSpecial(-4163)
$ws.Cells.Item(18, 3).Value = '와이바이오로직스'
$ws.Cells.Item(18, 4).Value = '유안타'
$ws.Cells.Item(18, 5).Value = '유안타'
$ws.Cells.Item(18, 6).Formula = '=""&"2023-11-28"'
$ws.Cells.Item(18, 6).Copy()
$ws.Cells.Item(18, 6).PasteSpecial(-4163)
$ws.Cells.Item(18, 7).Formula = '=""&"2023-12-05"'
$ws.Cells.Item(18, 7).Copy()
$ws.Cells.Item(18, 7).PasteSpecial(-4163)
$ws.Cells.Item(18, 8).Value = 13500
$ws.Cells.Item(18, 9).Value = 1500000
$ws.Cells.Item(18, 10).Value = 9000
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 100

# Row 19: 이베스트 / LS머트리얼즈
$ws.Cells.Item(19, 1).Value = '이베스트'
$ws.Cells.Item(19, 2).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(19, 2).PasteSpecial(-4163)
$ws.Cells.Item(19, 3).Value = 'LS머트리얼즈'
$ws.Cells.Item(19, 4).Value = '키움, KB'
$ws.Cells.Item(19, 5).Value = '키움, KB, 이베스트, 하이, NH'
$ws.Cells.Item(19, 6).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(19, 6).Copy()
$ws.Cells.Item(19, 6).PasteSpecial(-4163)
$ws.Cells.Item(19, 7).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(19, 7).Copy()
$ws.Cells.Item(19, 7).PasteSpecial(-4163)
$ws.Cells.Item(19, 8).Value = 6581.256
$ws.Cells.Item(19, 9).Value = 14625000
$ws.Cells.Item(19, 10).Value = 6000
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 7.5

# Row 20: 키움 / 블루엠텍
$ws.Cells.Item(20, 1).Value = '키움'
$ws.Cells.Item(20, 2).Formula = '=""&"2023-12-04"'
$ws.Cells.Item(20, 2).Copy()
$ws.Cells.Item(20, 2).PasteSpecial(-4163)
$ws.Cells.Item(20, 3).Value = '블루엠텍'
$ws.Cells.Item(20, 4).Value = '하나'
$ws.Cells.Item(20, 5).Value = '하나, 키움'
$ws.Cells.Item(20, 6).Formula = '=""&"2023-12-07"'
$ws.Cells.Item(20, 6).Copy()
$ws.Cells.Item(20, 6).PasteSpecial(-4163)
$ws.Cells.Item(20, 7).Formula = '=""&"2023-12-13"'
$ws.Cells.Item(20, 7).Copy()
$ws.Cells.Item(20, 7).PasteSpecial(-4163)
$ws.Cells.Item(20, 8).Value = 7980
$ws.Cells.Item(20, 9).Value = 1400000
$ws.Cells.Item(20, 10).Value = 19000
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 30

# Row 21: 키움 / LS머트리얼즈
$ws.Cells.Item(21, 1).Value = '키움'
$ws.Cells.Item(21, 2).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(21, 2).Copy()
$ws.Cells.Item(21, 2).PasteSpecial(-4163)
$ws.Cells.Item(21, 3).Value = 'LS머트리얼즈'
$ws.Cells.Item(21, 4).Value = '키움, KB'
$ws.Cells.Item(21, 5).Value = '키움, KB, 이베스트, 하이, NH'
$ws.Cells.Item(21, 6).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(21, 6).Copy()
$ws.Cells.Item(21, 6).PasteSpecial(-4163)
$ws.Cells.Item(21, 7).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(21, 7).Copy()
$ws.Cells.Item(21, 7).PasteSpecial(-4163)
$ws.Cells.Item(21, 8).Value = 36196.872
$ws.Cells.Item(21, 9).Value = 14625000
$ws.Cells.Item(21, 10).Value = 6000
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 41.25

# Row 22: 하나 / 블루엠텍
$ws.Cells.Item(22, 1).Value = '하나'
$ws.Cells.Item(22, 2).Formula = '=""&"2023-12-04"'
$ws.Cells.Item(22, 2).Copy()
$ws.Cells.Item(22, 2).PasteSpecial(-4163)
$ws.Cells.Item(22, 3).Value = '블루엠텍'
$ws.Cells.Item(22, 4).Value = '하나'
$ws.Cells.Item(22, 5).Value = '하나, 키움'
$ws.Cells.Item(22, 6).Formula = '=""&"2023-12-07"'
$ws.Cells.Item(22, 6).Copy()
$ws.Cells.Item(22, 6).PasteSpecial(-4163)
$ws.Cells.Item(22, 7).Formula = '=""&"2023-12-13"'
$ws.Cells.Item(22, 7).Copy()
$ws.Cells.Item(22, 7).PasteSpecial(-4163)
$ws.Cells.Item(22, 8).Value = 18620
$ws.Cells.Item(22, 9).Value = 1400000
$ws.Cells.Item(22, 10).Value = 19000
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 70

# Row 23: 하나 / 하나30호스팩
$ws.Cells.Item(23, 1).Value = '하나'
$ws.Cells.Item(23, 2).Formula = '=""&"2023-12-13"'
$ws.Cells.Item(23, 2).Copy()
$ws.Cells.Item(23, 2).PasteSpecial(-4163)
$ws.Cells.Item(23, 3).Value = '하나30호스팩'
$ws.Cells.Item(23, 4).Value = '하나'
$ws.Cells.Item(23, 5).Value = '하나'
$ws.Cells.Item(23, 6).Formula = '=""&"2023-12-18"'
$ws.Cells.Item(23, 6).Copy()
$ws.Cells.Item(23, 6).PasteSpecial(-4163)
$ws.Cells.Item(23, 7).Formula = '=""&"2023-12-22"'
$ws.Cells.Item(23, 7).Copy()
$ws.Cells.Item(23, 7).PasteSpecial(-4163)
$ws.Cells.Item(23, 8).Value = 14000
$ws.Cells.Item(23, 9).Value = 7000000
$ws.Cells.Item(23, 10).Value = 2000
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 100

# Row 24: 하나 / 에이텀
$ws.Cells.Item(24, 1).Value = '하나'
$ws.Cells.Item(24, 2).Formula = '=""&"2023-11-21"'
$ws.Cells.Item(24, 2).Copy()
$ws.Cells.Item(24, 2).PasteSpecial(-4163)
$ws.Cells.Item(24, 3).Value = '에이텀'
$ws.Cells.Item(24, 4).Value = '하나'
$ws.Cells.Item(24, 5).Value = '하나'
$ws.Cells.Item(24, 6).Formula = '=""&"2023-11-24"'
$ws.Cells.Item(24, 6).Copy()
$ws.Cells.Item(24, 6).PasteSpecial(-4163)
$ws.Cells.Item(24, 7).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(24, 7).Copy()
$ws.Cells.Item(24, 7).PasteSpecial(-4163)
$ws.Cells.Item(24, 8).Value = 11700
$ws.Cells.Item(24, 9).Value = 650000
$ws.Cells.Item(24, 10).Value = 18000
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 100

# Row 25: 하이 / 스톰테크
$ws.Cells.Item(25, 1).Value = '하이'
$ws.Cells.Item(25, 2).Formula = '=""&"2023-11-09"'
$ws.Cells.Item(25, 2).Copy()
$ws.Cells.Item(25, 2).PasteSpecial(-4163)
$ws.Cells.Item(25, 3).Value = '스톰테크'
$ws.Cells.Item(25, 4).Value = '하이'
$ws.Cells.Item(25, 5).Value = '하이'
$ws.Cells.Item(25, 6).Formula = '=""&"2023-11-14"'
$ws.Cells.Item(25, 6).Copy()
$ws.Cells.Item(25, 6).PasteSpecial(-4163)
$ws.Cells.Item(25, 7).Formula = '=""&"2023-11-20"'
$ws.Cells.Item(25, 7).Copy()
$ws.Cells.Item(25, 7).PasteSpecial(-4163)
$ws.Cells.Item(25, 8).Value = 36850
$ws.Cells.Item(25, 9).Value = 3350000
$ws.Cells.Item(25, 10).Value = 11000
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 100

# Row 26: 하이 / LS머트리얼즈
$ws.Cells.Item(26, 1).Value = '하이'
$ws.Cells.Item(26, 2).Formula = '=""&"2023-12-01"'
$ws.Cells.Item(26, 2).Copy()
$ws.Cells.Item(26, 2).PasteSpecial(-4163)
$ws.Cells.Item(26, 3).Value = 'LS머트리얼즈'
$ws.Cells.Item(26, 4).Value = '키움, KB'
$ws.Cells.Item(26, 5).Value = '키움, KB, 이베스트, 하이, NH'
$ws.Cells.Item(26, 6).Formula = '=""&"2023-12-06"'
$ws.Cells.Item(26, 6).Copy()
$ws.Cells.Item(26, 6).PasteSpecial(-4163)
$ws.Cells.Item(26, 7).Formula = '=""&"2023-12-12"'
$ws.Cells.Item(26, 7).Copy()
$ws.Cells.Item(26, 7).PasteSpecial(-4163)
$ws.Cells.Item(26, 8).Value = 4387.5
$ws.Cells.Item(26, 9).Value = 14625000
$ws.Cells.Item(26, 10).Value = 6000
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 5

# Row 27: 하이 / 에코프로머티
$ws.Cells.Item(27, 1).Value = '하이'
$ws.Cells.Item(27, 2).Formula = '=""&"2023-11-08"'
$ws.Cells.Item(27, 2).Copy()
$ws.Cells.Item(27, 2).PasteSpecial(-4163)
$ws.Cells.Item(27, 3).Value = '에코프로머티'
$ws.Cells.Item(27, 4).Value = '미래'
$ws.Cells.Item(27, 5).Value = '미래, NH, 하이'
$ws.Cells.Item(27, 6).Formula = '=""&"2023-11-13"'
$ws.Cells.Item(27, 6).Copy()
$ws.Cells.Item(27, 6).PasteSpecial(-4163)
$ws.Cells.Item(27, 7).Formula = '=""&"2023-11-17"'
$ws.Cells.Item(27, 7).Copy()
$ws.Cells.Item(27, 7).PasteSpecial(-4163)
$ws.Cells.Item(27, 8).Value = 12576.7488
$ws.Cells.Item(27, 9).Value = 11580800
$ws.Cells.Item(27, 10).Value = 36200
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 3

$excel.CutCopyMode = 0
